{"js": "// Replace each two-digit multiplication expression with its new value.\n// Pairs are unique literal strings in the document, so Body.search finds\n// exactly one match per pair and we can safely replace in place.\nconst replacements = [\n  [\"14\u00d719=266\", \"28\u00d725=700\"],\n  [\"57\u00d788=5016\", \"23\u00d778=1794\"],\n  [\"86\u00d734=2924\", \"50\u00d721=1050\"],\n  [\"12\u00d750=600\", \"89\u00d753=4717\"],\n  [\"73\u00d728=2044\", \"12\u00d783=996\"],\n  [\"90\u00d750=4500\", \"13\u00d798=1274\"],\n  [\"67\u00d795=6365\", \"67\u00d770=4690\"],\n  [\"28\u00d736=1008\", \"22\u00d742=924\"],\n  [\"69\u00d790=6210\", \"83\u00d749=4067\"],\n  [\"47\u00d733=1551\", \"79\u00d785=6715\"],\n  [\"74\u00d759=4366\", \"21\u00d793=1953\"],\n  [\"81\u00d727=2187\", \"58\u00d793=5394\"],\n  [\"21\u00d774=1554\", \"50\u00d713=650\"],\n  [\"70\u00d782=5740\", \"24\u00d791=2184\"],\n  [\"92\u00d735=3220\", \"98\u00d738=3724\"],\n  [\"70\u00d796=6720\", \"47\u00d726=1222\"],\n  [\"66\u00d726=1716\", \"20\u00d767=1340\"],\n  [\"90\u00d733=2970\", \"62\u00d735=2170\"],\n  [\"22\u00d781=1782\", \"13\u00d768=884\"],\n  [\"72\u00d790=6480\", \"89\u00d758=5162\"],\n  [\"52\u00d711=572\", \"75\u00d745=3375\"],\n  [\"72\u00d731=2232\", \"75\u00d714=1050\"],\n  [\"38\u00d790=3420\", \"27\u00d744=1188\"],\n  [\"58\u00d776=4408\", \"63\u00d735=2205\"],\n  [\"44\u00d715=660\", \"82\u00d715=1230\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each left-hand expression occurs exactly once in the document, so\n# Find/Replace with wdReplaceAll is equivalent to a single targeted swap.\n$replacements = @(\n    @{ Old = \"14\u00d719=266\"; New = \"28\u00d725=700\" },\n    @{ Old = \"57\u00d788=5016\"; New = \"23\u00d778=1794\" },\n    @{ Old = \"86\u00d734=2924\"; New = \"50\u00d721=1050\" },\n    @{ Old = \"12\u00d750=600\"; New = \"89\u00d753=4717\" },\n    @{ Old = \"73\u00d728=2044\"; New = \"12\u00d783=996\" },\n    @{ Old = \"90\u00d750=4500\"; New = \"13\u00d798=1274\" },\n    @{ Old = \"67\u00d795=6365\"; New = \"67\u00d770=4690\" },\n    @{ Old = \"28\u00d736=1008\"; New = \"22\u00d742=924\" },\n    @{ Old = \"69\u00d790=6210\"; New = \"83\u00d749=4067\" },\n    @{ Old = \"47\u00d733=1551\"; New = \"79\u00d785=6715\" },\n    @{ Old = \"74\u00d759=4366\"; New = \"21\u00d793=1953\" },\n    @{ Old = \"81\u00d727=2187\"; New = \"58\u00d793=5394\" },\n    @{ Old = \"21\u00d774=1554\"; New = \"50\u00d713=650\" },\n    @{ Old = \"70\u00d782=5740\"; New = \"24\u00d791=2184\" },\n    @{ Old = \"92\u00d735=3220\"; New = \"98\u00d738=3724\" },\n    @{ Old = \"70\u00d796=6720\"; New = \"47\u00d726=1222\" },\n    @{ Old = \"66\u00d726=1716\"; New = \"20\u00d767=1340\" },\n    @{ Old = \"90\u00d733=2970\"; New = \"62\u00d735=2170\" },\n    @{ Old = \"22\u00d781=1782\"; New = \"13\u00d768=884\" },\n    @{ Old = \"72\u00d790=6480\"; New = \"89\u00d758=5162\" },\n    @{ Old = \"52\u00d711=572\"; New = \"75\u00d745=3375\" },\n    @{ Old = \"72\u00d731=2232\"; New = \"75\u00d714=1050\" },\n    @{ Old = \"38\u00d790=3420\"; New = \"27\u00d744=1188\" },\n    @{ Old = \"58\u00d776=4408\"; New = \"63\u00d735=2205\" },\n    @{ Old = \"44\u00d715=660\"; New = \"82\u00d715=1230\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $found = $find.Execute(\n        $pair.Old,        # FindText\n        $true,            # MatchCase\n        $false,           # MatchWholeWord\n        $false,           # MatchWildcards\n        $false,           # MatchSoundsLike\n        $false,           # MatchAllWordForms\n        $true,            # Forward\n        \"wdFindContinue\", # Wrap\n        $false,           # Format\n        $pair.New,        # ReplaceWith\n        \"wdReplaceAll\"    # Replace\n    )\n    if (-not $found) {\n        throw \"No match found for $($pair.Old)\"\n    }\n}\n\n"}
